$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 296 (shifts old rows 296..415 down to 297..416,
# carrying all their original formatting/values with them).
$ws.Rows.Item(296).Insert()

# Populate the newly-inserted row 296 with a new daily record for
# "Feria Lagunitas de Puerto Montt" / Pepino ensalada, identical to the
# entry that is now in row 297 except for a new date (2023-04-11 = serial 45027).
$ws.Range("A296").Value = 4
$ws.Range("B296").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C296").Value = "Los Lagos"
$ws.Range("D296").Value = 45027
$ws.Range("E296").Value = 10
$ws.Range("F296").Value = 100112043
$ws.Range("G296").Value = "Pepino ensalada"
$ws.Range("H296").Value = "Sin especificar"
$ws.Range("I296").Value = "Primera"
$ws.Range("J296").Value = 400
$ws.Range("K296").Value = 13000
$ws.Range("L296").Value = 13000
$ws.Range("M296").Value = 13000
$ws.Range("N296").Value = "$/caja 60 unidades"
$ws.Range("O296").Value = "Región de Arica y Parinacota"
$ws.Range("P296").Value = 217
$ws.Range("Q296").Value = 60
$ws.Range("R296").Value = "Hortaliza"
